$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 5 with data (weekly update for Espárragos)
$ws.Cells.Item(5, 1).Value = 7
$ws.Cells.Item(5, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(5, 3).Value = "Ñuble"
$ws.Cells.Item(5, 4).Value = (Get-Date -Year 2021 -Month 11 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(5, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(5, 5).Value = 16
$ws.Cells.Item(5, 6).Value = 300000000
$ws.Cells.Item(5, 7).Value = "Espárragos"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 500
$ws.Cells.Item(5, 11).Value = 900
$ws.Cells.Item(5, 12).Value = 1000
$ws.Cells.Item(5, 13).Value = 950
$ws.Cells.Item(5, 14).Value = "`$/kilo"
$ws.Cells.Item(5, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(5, 16).Value = 950
$ws.Cells.Item(5, 17).Value = 1
$ws.Cells.Item(5, 18).Value = "Hortaliza"
